$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 2001.0869
$ws.Range("J32").Value = 2196.3333
$ws.Range("L32").Value = 2196.3333
$ws.Range("N32").Value = -2848.3333
# Row 88
$ws.Range("H88").Value = 3111
$ws.Range("J88").Value = 3166.5
$ws.Range("L88").Value = 3166.5
$ws.Range("N88").Value = -3978.5
# Row 91
$ws.Range("H91").Value = 3111
$ws.Range("J91").Value = 3166.5
$ws.Range("L91").Value = 3166.5
$ws.Range("N91").Value = -5974.5
# Row 111
$ws.Range("H111").Value = 2564.1904
$ws.Range("I111").Value = 2268.2144
$ws.Range("J111").Value = 3156.1428
$ws.Range("K111").Value = 6804.6432
$ws.Range("L111").Value = 9468.428400000001
$ws.Range("M111").Value = -3737.6432
$ws.Range("N111").Value = -15602.4284
# Row 116
$ws.Range("H116").Value = 12827.117
$ws.Range("J116").Value = 11938.896
$ws.Range("L116").Value = 11938.896
$ws.Range("N116").Value = -18822.896
# Row 137
$ws.Range("H137").Value = 1894
$ws.Range("I137").Value = 893.6667
$ws.Range("K137").Value = 2681.0001
$ws.Range("M137").Value = -131.0001000000002
# Row 138
$ws.Range("H138").Value = 2924.5789
$ws.Range("J138").Value = 3560.9644
$ws.Range("L138").Value = 10682.8932
$ws.Range("N138").Value = -20962.8932

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 7900
$ws.Range("J8").Value = 800
$ws.Range("L8").Value = 800
$ws.Range("N8").Value = -1088
# Row 13
$ws.Range("H13").Value = 8424.5
$ws.Range("J13").Value = 11166
$ws.Range("L13").Value = 11166
$ws.Range("N13").Value = -11454
# Row 61
$ws.Range("H61").Value = 5745.854
$ws.Range("I61").Value = 4418.8623
$ws.Range("K61").Value = 4418.8623
$ws.Range("M61").Value = -4206.8623
# Row 63
$ws.Range("H63").Value = 2926.1667
$ws.Range("I63").Value = 2926.1667
$ws.Range("K63").Value = 2926.1667
$ws.Range("M63").Value = -2240.1667
# Row 66
$ws.Range("H66").Value = 2926.1667
$ws.Range("I66").Value = 2926.1667
$ws.Range("K66").Value = 14630.8335
$ws.Range("M66").Value = -11198.8335
# Row 136
$ws.Range("H136").Value = 5745.854
$ws.Range("I136").Value = 4418.8623
$ws.Range("K136").Value = 13256.5869
$ws.Range("M136").Value = -10706.5869

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 11962.637
$ws.Range("I86").Value = 3012.8572
$ws.Range("K86").Value = 3012.8572
$ws.Range("M86").Value = -1889.8572
# Row 89
$ws.Range("H89").Value = 11962.637
$ws.Range("I89").Value = 3012.8572
$ws.Range("K89").Value = 15064.286
$ws.Range("M89").Value = -9448.286
# Row 94
$ws.Range("H94").Value = 2707750
$ws.Range("I94").Value = 3032778
$ws.Range("J94").Value = 26269
$ws.Range("K94").Value = 3032778
$ws.Range("L94").Value = 26269
$ws.Range("M94").Value = -3032327
$ws.Range("N94").Value = -27171
# Row 107
$ws.Range("H107").Value = 1793.3182
$ws.Range("I107").Value = 1753.3334
$ws.Range("K107").Value = 1753.3334
$ws.Range("M107").Value = 166.6666
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
# Row 134
$ws.Range("H134").Value = 4918.6807
$ws.Range("I134").Value = 4355.5137
$ws.Range("K134").Value = 13066.5411
$ws.Range("M134").Value = -10531.5411

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3305.6316
$ws.Range("I31").Value = 1801
$ws.Range("J31").Value = 5374.5
$ws.Range("K31").Value = 1801
$ws.Range("L31").Value = 5374.5
$ws.Range("M31").Value = -1506
$ws.Range("N31").Value = -5964.5
# Row 34
$ws.Range("H34").Value = 3305.6316
$ws.Range("I34").Value = 1801
$ws.Range("J34").Value = 5374.5
$ws.Range("K34").Value = 1801
$ws.Range("L34").Value = 5374.5
$ws.Range("M34").Value = -1599
$ws.Range("N34").Value = -5778.5
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 86
$ws.Range("H86").Value = 5212
$ws.Range("I86").Value = 4949.6665
$ws.Range("J86").Value = 5999
$ws.Range("K86").Value = 4949.6665
$ws.Range("L86").Value = 5999
$ws.Range("M86").Value = -3826.6665
$ws.Range("N86").Value = -8245
# Row 89
$ws.Range("H89").Value = 5212
$ws.Range("I89").Value = 4949.6665
$ws.Range("J89").Value = 5999
$ws.Range("K89").Value = 24748.3325
$ws.Range("L89").Value = 29995
$ws.Range("M89").Value = -19132.3325
$ws.Range("N89").Value = -41227
# Row 107
$ws.Range("H107").Value = 1258.7391
$ws.Range("I107").Value = 1111.625
$ws.Range("K107").Value = 1111.625
$ws.Range("M107").Value = 808.375
# Row 141
$ws.Range("H141").Value = 77220
$ws.Range("J141").Value = 77220
$ws.Range("L141").Value = 77220
$ws.Range("N141").Value = -87580

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 112
$ws.Range("H112").Value = 6842.8335
$ws.Range("J112").Value = 8019.6665
$ws.Range("L112").Value = 24058.9995
$ws.Range("N112").Value = -26274.9995
# Row 128
$ws.Range("H128").Value = 224950
$ws.Range("I128").Value = 224950
$ws.Range("K128").Value = 674850
$ws.Range("M128").Value = -669870
# Row 131
$ws.Range("H131").Value = 24693674
$ws.Range("J131").Value = 83336570
$ws.Range("L131").Value = 250009710
$ws.Range("N131").Value = -250019790

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 25604736
$ws.Range("J11").Value = 82715140
$ws.Range("L11").Value = 82715140
$ws.Range("N11").Value = -82715418
# Row 25
$ws.Range("H25").Value = 9999.5
$ws.Range("J25").Value = 9999.5
$ws.Range("L25").Value = 9999.5
$ws.Range("N25").Value = -11057.5
# Row 80
$ws.Range("H80").Value = 2472.9524
$ws.Range("I80").Value = 2361.3333
$ws.Range("J80").Value = 2556.6667
$ws.Range("K80").Value = 2361.3333
$ws.Range("L80").Value = 2556.6667
$ws.Range("M80").Value = -1363.3333
$ws.Range("N80").Value = -4552.6667
# Row 83
$ws.Range("H83").Value = 2472.9524
$ws.Range("I83").Value = 2361.3333
$ws.Range("J83").Value = 2556.6667
$ws.Range("K83").Value = 11806.6665
$ws.Range("L83").Value = 12783.3335
$ws.Range("M83").Value = -6814.666499999999
$ws.Range("N83").Value = -22767.3335
# Row 97
$ws.Range("H97").Value = 1327.6154
$ws.Range("I97").Value = 613.4545000000001
$ws.Range("K97").Value = 613.4545000000001
$ws.Range("M97").Value = -117.4545000000001
# Row 113
$ws.Range("H113").Value = 2388
$ws.Range("I113").Value = 2000.1111
$ws.Range("J113").Value = 3086.2
$ws.Range("K113").Value = 2000.1111
$ws.Range("L113").Value = 3086.2
$ws.Range("M113").Value = 169.8888999999999
$ws.Range("N113").Value = -7426.2
# Row 135
$ws.Range("H135").Value = 92000
$ws.Range("J135").Value = 92000
$ws.Range("L135").Value = 92000
$ws.Range("N135").Value = -102140
# Row 138
$ws.Range("H138").Value = 70000
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2031.4482
$ws.Range("I68").Value = 1842.1765
$ws.Range("J68").Value = 2299.5833
$ws.Range("K68").Value = 1842.1765
$ws.Range("L68").Value = 2299.5833
$ws.Range("M68").Value = -1093.1765
$ws.Range("N68").Value = -3797.5833
# Row 71
$ws.Range("H71").Value = 2031.4482
$ws.Range("I71").Value = 1842.1765
$ws.Range("J71").Value = 2299.5833
$ws.Range("K71").Value = 9210.8825
$ws.Range("L71").Value = 11497.9165
$ws.Range("M71").Value = -5466.8825
$ws.Range("N71").Value = -18985.9165
# Row 132
$ws.Range("H132").Value = 3370.1482
$ws.Range("I132").Value = 2643.1365
$ws.Range("K132").Value = 7929.4095
$ws.Range("M132").Value = -5399.4095
# Row 136
$ws.Range("H136").Value = 2562.2
$ws.Range("I136").Value = 1642
$ws.Range("K136").Value = 4926
$ws.Range("M136").Value = -2376
# Row 140
$ws.Range("H140").Value = 63498.6
$ws.Range("J140").Value = 64374.25
$ws.Range("L140").Value = 64374.25
$ws.Range("N140").Value = -74734.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 49
$ws.Range("H49").Value = 36582.668
$ws.Range("J49").Value = 49999
$ws.Range("L49").Value = 49999
$ws.Range("N49").Value = -50459
# Row 132
$ws.Range("H132").Value = 3313.3242
$ws.Range("I132").Value = 3080.4583
$ws.Range("J132").Value = 3743.2307
$ws.Range("K132").Value = 9241.374899999999
$ws.Range("L132").Value = 11229.6921
$ws.Range("M132").Value = -6711.374899999999
$ws.Range("N132").Value = -16289.6921
# Row 136
$ws.Range("H136").Value = 6028.3105
$ws.Range("I136").Value = 6028.3105
$ws.Range("K136").Value = 18084.9315
$ws.Range("M136").Value = -15534.9315
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
